$wb = $excel.ActiveWorkbook
Write-Host $wb.Windows.Count
$win = $wb.Windows.Item(1)
Write-Host ($win | Get-Member | Out-String)
